# The Act (previous commit) had assumed contemporaneous parameters.
# new_medicare_levy = medicare_levy
#
# Update the "indiv" sheet's historical medicare-levy parameters for
# 2003-04 .. 2006-07 (rows 2-7): correct lower/upper thresholds and the
# taper rate, and backfill the family-threshold columns (J/K/L) that
# were already present for later years.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("indiv")

$fmt = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# Row 2: 2003-04, individual
$ws.Range("F2").Value = 15529
$ws.Range("G2").Value = 16789
$ws.Range("H2").Value = 0.2
$ws.Range("J2").Value = 26205
$ws.Range("K2").Value = 26206
$ws.Range("L2").Value = 2406
$ws.Range("J8").Copy() | Out-Null
$ws.Range("L2").PasteSpecial($fmt) | Out-Null

# Row 3: 2004-05, individual
$ws.Range("F3").Value = 15902
$ws.Range("G3").Value = 17192
$ws.Range("H3").Value = 0.2
$ws.Range("J3").Value = 26834
$ws.Range("K3").Value = 26835
$ws.Range("L3").Value = 2464
$ws.Range("J8").Copy() | Out-Null
$ws.Range("L3").PasteSpecial($fmt) | Out-Null

# Row 4: 2005-06, individual
$ws.Range("F4").Value = 16284
$ws.Range("G4").Value = 17605
$ws.Range("H4").Value = 0.2
$ws.Range("J4").Value = 27478
$ws.Range("K4").Value = 27479
$ws.Range("L4").Value = 2523
$ws.Range("J8").Copy() | Out-Null
$ws.Range("L4").PasteSpecial($fmt) | Out-Null

# Row 5: 2006-07, individual (taper unchanged at 0.1)
$ws.Range("F5").Value = 16740
$ws.Range("G5").Value = 19695
$ws.Range("J5").Value = 28247
$ws.Range("K5").Value = 28248
$ws.Range("L5").Value = 2594
$ws.Range("J8").Copy() | Out-Null
$ws.Range("L5").PasteSpecial($fmt) | Out-Null

# Row 6: 2006-07, sato/pto/sapto (no threshold changes, only new family cols)
$ws.Range("J6").Value = 33500
$ws.Range("K6").Value = 33501
$ws.Range("L6").Value = 2594
$ws.Range("J8").Copy() | Out-Null
$ws.Range("L6").PasteSpecial($fmt) | Out-Null

# Row 7: 2006-07, pto only (no threshold changes, only new family cols)
$ws.Range("J7").Value = 33500
$ws.Range("K7").Value = 33501
$ws.Range("L7").Value = 2594
$ws.Range("J8").Copy() | Out-Null
$ws.Range("L7").PasteSpecial($fmt) | Out-Null

$excel.CutCopyMode = $false

# Restore the sheet's last active-cell selection.
$ws.Range("L3").Select() | Out-Null
